$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row: _old -> _FV2310, _new -> _FV2404 ---
$oldHeaders = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 1          # columns A..J
    $ws.Cells.Item(1, $col).Value = $oldHeaders[$i] + "_FV2310"
}

# Column K (11) stays "diff"

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $col = $i + 12          # columns L..U
    $ws.Cells.Item(1, $col).Value = $oldHeaders[$i] + "_FV2404"
}

# --- 2. Freeze the header row (pane split after row 1) ---
$ws.Range("A2").Select() | Out-Null
[void]($excel.ActiveWindow.FreezePanes = $true)

# --- 3. Turn the used range into an Excel Table (ListObject) ---
$rng = $ws.Range("A1:U65")
$tbl = $ws.ListObjects.Add(1, $rng, $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
